$d = $word.ActiveDocument

function Get-ParagraphIndexContaining($doc, $needle) {
    for ($i = 1; $i -le $doc.Paragraphs.Count; $i++) {
        $p = $doc.Paragraphs.Item($i)
        if ($p.Range.Text.Contains($needle)) {
            return $i
        }
    }
    return -1
}

# ------------------------------------------------------------------
# 1) Remove the paragraph holding the italic list of key terms that
#    used to echo the "O" heading paragraph's contents.
# ------------------------------------------------------------------
$idx = Get-ParagraphIndexContaining $d "Obstiné, Offrandes céréalières, Oint, Olivier, Onésime, Osée, Otniel"
if ($idx -gt 0) {
    $d.Paragraphs.Item($idx).Range.Delete()
}

# ------------------------------------------------------------------
# 2) Remove the "This PDF version is provided under the same
#    license." paragraph -- it is merged away.
# ------------------------------------------------------------------
$idx = Get-ParagraphIndexContaining $d "This PDF version is provided under the same license."
if ($idx -gt 0) {
    $d.Paragraphs.Item($idx).Range.Delete()
}

# ------------------------------------------------------------------
# 3) Rewrite the license-description paragraph (the one that starts
#    with a bold run) with the new resource-data text.
# ------------------------------------------------------------------
$idx = Get-ParagraphIndexContaining $d "is based on"
$p = $d.Paragraphs.Item($idx)
$r = $p.Range
$r.MoveEnd(1, -1) | Out-Null

$t1 = "Biblica Study Notes (Key Terms)"
$t2 = " © 2023 Biblica Inc. Released under CC BY-SA 4.0 license. "
$t3 = "Biblica Study Notes"
$t4 = " has been adapted in the following languages: Tok Pisin, Arabic (عربي), French (Français), Hindi (हिंदी), Indonesian (Bahasa Indonesia), Portuguese (Português), Russian (Русский), Spanish (Español), Swahili (Kiswahili), and Simplified Chinese (简体中文)from Biblica Study Notes © 2023 Biblica Inc. Released under CC BY-SA 4.0 license by Mission Mutual."

$full = $t1 + $t2 + $t3 + $t4
$r.Text = $full

$start = $r.Start
$pos2 = $start + $t1.Length
$pos3 = $pos2 + $t2.Length
$pos4 = $pos3 + $t3.Length
$pos5 = $pos4 + $t4.Length

# The first run inherits the original Bold formatting of the
# paragraph's lead run; explicitly un-bold the remaining three
# pieces so they serialize as distinct, correctly-tagged runs.
$d.Range($pos2, $pos3).Font.Bold = 0
$d.Range($pos3, $pos4).Font.Bold = 0
$d.Range($pos4, $pos5).Font.Bold = 0

# ------------------------------------------------------------------
# 4) Remove the "License Information" Heading 2 paragraph entirely.
# ------------------------------------------------------------------
$idx = Get-ParagraphIndexContaining $d "License Information"
if ($idx -gt 0) {
    $d.Paragraphs.Item($idx).Range.Delete()
}

Write-Host "Done. Paragraph count: $($d.Paragraphs.Count)"
